# Generation Data.xlsx -- debugging edit: all indexes to 0, all variables
# created for all nodes.
#
# Sheet "Generation_investor" (sheet1) and "Generation_rival" (sheet2) are
# both expanded from 6 data rows to 24 data rows (one row per node, 1..24),
# inserting the missing "all zero" placeholder rows that were previously
# skipped, and re-numbering column B (Node) as a running count using
# formulas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Generation_investor
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Generation_investor")

# Existing (pre-edit) data, keyed by the Node number (column B value):
#   Node 7  -> A=3  C=350 D=20.7
#   Node 13 -> A=4  C=591 D=20.93
#   Node 15 -> A=5  C=215 (formula) D=18.315 (formula) -- handled separately below
$investorData = @{
    7  = @{ A = 3; C = 350; D = 20.7 }
    13 = @{ A = 4; C = 591; D = 20.93 }
}

for ($node = 3; $node -le 24; $node++) {
    $r = $node + 1   # row 4 holds node 3, row 5 holds node 4, ... row r = node + 1

    if ($investorData.ContainsKey($node)) {
        $d = $investorData[$node]
        $ws1.Cells.Item($r, 1).Value = $d.A
        $ws1.Cells.Item($r, 3).Value = $d.C
        $ws1.Cells.Item($r, 4).Value = $d.D
    } else {
        $ws1.Cells.Item($r, 1).Value = 0
        $ws1.Cells.Item($r, 3).Value = 0
        $ws1.Cells.Item($r, 4).Value = 0
    }
}

# Row 16 (node 15) keeps its original formulas but moves down from row 6.
$ws1.Cells.Item(16, 1).Value = 5
$ws1.Cells.Item(16, 3).Formula = "=60+155"
$ws1.Cells.Item(16, 4).Formula = "=(26.11+10.52)/2"

# Column B (Node index) becomes a running count via formulas, rows 4-25.
$ws1.Cells.Item(4, 2).Value = 3
for ($r = 5; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 2).Formula = "=B" + ($r - 1) + "+1"
}

# ---------------------------------------------------------------------
# Sheet 2: Generation_rival
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Generation_rival")

$rivalData = @{
    16 = @{ A = 7;  C = 155; D = 10.52 }
    18 = @{ A = 8;  C = 400; D = 6.02 }
    21 = @{ A = 9;  C = 400; D = 5.47 }
    22 = @{ A = 10; C = 300; D = 0 }
}

for ($node = 1; $node -le 23; $node++) {
    $r = $node + 1

    if ($rivalData.ContainsKey($node)) {
        $d = $rivalData[$node]
        $ws2.Cells.Item($r, 1).Value = $d.A
        $ws2.Cells.Item($r, 3).Value = $d.C
        $ws2.Cells.Item($r, 4).Value = $d.D
    } else {
        $ws2.Cells.Item($r, 1).Value = 0
        $ws2.Cells.Item($r, 3).Value = 0
        $ws2.Cells.Item($r, 4).Value = 0
    }
}

# Row 24 (node 23) keeps its original formulas, moved down from row 6.
$ws2.Cells.Item(24, 1).Value = 11
$ws2.Cells.Item(24, 3).Formula = "=350+310"
$ws2.Cells.Item(24, 4).Formula = "=(10.52+10.89)/2"

# Row 25 is a brand-new all-zero row for node 24.
$ws2.Cells.Item(25, 1).Value = 0
$ws2.Cells.Item(25, 2).Value = 24
$ws2.Cells.Item(25, 3).Value = 0
$ws2.Cells.Item(25, 4).Value = 0

# Column B (Node index) becomes a running count via formulas, rows 2-24.
$ws2.Cells.Item(2, 2).Value = 1
for ($r = 3; $r -le 24; $r++) {
    $ws2.Cells.Item($r, 2).Formula = "=B" + ($r - 1) + "+1"
}

$ws2.Range("A1:D25").EntireColumn.AutoFit() | Out-Null

$wb.Save()
